# Insert a new weekly price record at row 301 for
# "Hortaliza, Terminal La Palmera de La Serena - Cebollín".
# This pushes the existing rows 301..363 down to 302..364
# (dimension grows from A1:R363 to A1:R364).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 301..363 down by one, leaving a blank row 301.
$ws.Rows.Item(301).Insert()

# Populate the newly inserted row 301 with the new record's values.
$ws.Cells.Item(301, 1).Value  = 8
$ws.Cells.Item(301, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(301, 3).Value  = "Coquimbo"
$ws.Cells.Item(301, 4).Value  = 45173
$ws.Cells.Item(301, 5).Value  = 4
$ws.Cells.Item(301, 6).Value  = 100112037
$ws.Cells.Item(301, 7).Value  = "Cebollín"
$ws.Cells.Item(301, 8).Value  = "Sin especificar"
$ws.Cells.Item(301, 9).Value  = "Primera"
$ws.Cells.Item(301, 10).Value = 1600
$ws.Cells.Item(301, 11).Value = 1000
$ws.Cells.Item(301, 12).Value = 1200
$ws.Cells.Item(301, 13).Value = 1100
$ws.Cells.Item(301, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(301, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(301, 16).Value = 183
$ws.Cells.Item(301, 17).Value = 6
$ws.Cells.Item(301, 18).Value = "Hortaliza"
